$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 928, shifting existing rows 928-994 down to 930-996
$ws.Rows("928:929").Insert()

# New row 928 data
$ws.Cells.Item(928, 1).Value = 5
$ws.Cells.Item(928, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(928, 3).Value = "Maule"
$ws.Cells.Item(928, 4).Value = 45021
$ws.Cells.Item(928, 5).Value = 7
$ws.Cells.Item(928, 6).Value = 100112020
$ws.Cells.Item(928, 7).Value = "Tomate"
$ws.Cells.Item(928, 8).Value = "Larga vida"
$ws.Cells.Item(928, 9).Value = "Primera"
$ws.Cells.Item(928, 10).Value = 2500
$ws.Cells.Item(928, 11).Value = 9000
$ws.Cells.Item(928, 12).Value = 9000
$ws.Cells.Item(928, 13).Value = 9000
$ws.Cells.Item(928, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(928, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(928, 16).Value = 500
$ws.Cells.Item(928, 17).Value = 18
$ws.Cells.Item(928, 18).Value = "Hortaliza"

# New row 929 data
$ws.Cells.Item(929, 1).Value = 5
$ws.Cells.Item(929, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(929, 3).Value = "Maule"
$ws.Cells.Item(929, 4).Value = 45021
$ws.Cells.Item(929, 5).Value = 7
$ws.Cells.Item(929, 6).Value = 100112020
$ws.Cells.Item(929, 7).Value = "Tomate"
$ws.Cells.Item(929, 8).Value = "Larga vida"
$ws.Cells.Item(929, 9).Value = "Primera"
$ws.Cells.Item(929, 10).Value = 1500
$ws.Cells.Item(929, 11).Value = 4500
$ws.Cells.Item(929, 12).Value = 4500
$ws.Cells.Item(929, 13).Value = 4500
$ws.Cells.Item(929, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(929, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(929, 16).Value = 450
$ws.Cells.Item(929, 17).Value = 10
$ws.Cells.Item(929, 18).Value = "Hortaliza"

# Ensure date formatting for column D on the new rows matches the rest of the column (style s="2")
$ws.Cells.Item(928, 4).NumberFormat = $ws.Cells.Item(930, 4).NumberFormat
$ws.Cells.Item(929, 4).NumberFormat = $ws.Cells.Item(930, 4).NumberFormat
